# Regenerate the localization-status report: the zh-cn / de-de handoff
# items have moved on from "Ready for handoff" to "In Translation", so
# update every cell that shows that status across the Overview, zh-cn and
# de-de sheets, then re-tighten the now-narrower status columns.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: zh-cn (col E) / de-de (col F) status columns ------
$ovw = $wb.Worksheets.Item(1)
$ovw.Range("E2").Value = $newStatus
$ovw.Range("F2").Value = $newStatus
$ovw.Range("E3").Value = $newStatus
$ovw.Range("F3").Value = $newStatus
$ovw.Range("E4").Value = $newStatus
$ovw.Range("F4").Value = $newStatus

# --- zh-cn sheet: Status column (col C) ---------------------------------
$zhcn = $wb.Worksheets.Item(2)
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("C4").Value = $newStatus

# --- de-de sheet: Status column (col C) ----------------------------------
$dede = $wb.Worksheets.Item(3)
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus
$dede.Range("C4").Value = $newStatus

# --- Shrink the status columns now that the text is shorter -------------
# ("In Translation" is shorter than "Ready for handoff", so the report
# generator re-sizes the columns that hold the status value.)
$ovw.Range("E1:F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5

Write-Output "Updated status to 'In Translation' and resized status columns."
